$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.033007874342867
$ws.Range("D2").Value = 1.042261326451732
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.050770146796068
$ws.Range("I2").Value = 1.038076159013012
$ws.Range("J2").Value = 1.038134962944198
$ws.Range("K2").Value = 1.045038513414612
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.05352350201329
$ws.Range("N2").Value = 1.039609233667323

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.033978056411331
$ws.Range("D3").Value = 1.043025207178789
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.051666640900114
$ws.Range("I3").Value = 1.038287010373691
$ws.Range("J3").Value = 1.038747573574859
$ws.Range("K3").Value = 1.045613312674327
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.054232292614091
$ws.Range("N3").Value = 1.040222714275347

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.034605859884538
$ws.Range("D4").Value = 1.043519167376833
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.052246729558984
$ws.Range("I4").Value = 1.038421650388943
$ws.Range("J4").Value = 1.039143405689302
$ws.Range("K4").Value = 1.045984264417345
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.05469029053768
$ws.Range("N4").Value = 1.040619108516795

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034869795255352
$ws.Range("D5").Value = 1.043726750062387
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.052490596956829
$ws.Range("I5").Value = 1.038477822747752
$ws.Range("J5").Value = 1.039309676812489
$ws.Range("K5").Value = 1.046139976583955
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.054882679358111
$ws.Range("N5").Value = 1.040785615764049

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034914111546709
$ws.Range("D6").Value = 1.043761599490119
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.052531543222156
$ws.Range("I6").Value = 1.038487229083293
$ws.Range("J6").Value = 1.039337586428232
$ws.Range("K6").Value = 1.046166107479666
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.054914973256422
$ws.Range("N6").Value = 1.040813565014648

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.034609386575934
$ws.Range("D7").Value = 1.043521941415516
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.05224998813593
$ws.Range("I7").Value = 1.03842240265831
$ws.Range("J7").Value = 1.039145627950249
$ws.Range("K7").Value = 1.045986345979187
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.054692861850975
$ws.Range("N7").Value = 1.040621333933608

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033335745188295
$ws.Range("D8").Value = 1.042519549435545
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.051073121094909
$ws.Range("I8").Value = 1.038147788508788
$ws.Range("J8").Value = 1.038342114728187
$ws.Range("K8").Value = 1.045232972268988
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.05376317240364
$ws.Range("N8").Value = 1.039816679630606

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031091687440191
$ws.Range("D9").Value = 1.040750793573975
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.048999358995015
$ws.Range("I9").Value = 1.037650162408478
$ws.Range("J9").Value = 1.036921904101659
$ws.Range("K9").Value = 1.043897951139892
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.052120106264003
$ws.Range("N9").Value = 1.038394452142088

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029595855895641
$ws.Range("D10").Value = 1.039570067261608
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.047616936166047
$ws.Range("I10").Value = 1.037309220526369
$ws.Range("J10").Value = 1.035972235484293
$ws.Range("K10").Value = 1.043002958712601
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.051021535871792
$ws.Range("N10").Value = 1.037443434886357

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028948199199532
$ws.Range("D11").Value = 1.039058444237506
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.047018365406286
$ws.Range("I11").Value = 1.037159415722592
$ws.Range("J11").Value = 1.035560347963875
$ws.Range("K11").Value = 1.042614246539351
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.050545096972295
$ws.Range("N11").Value = 1.037030962438416

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02870763828059
$ws.Range("D12").Value = 1.038868351492736
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.04679603470462
$ws.Range("I12").Value = 1.037103445312742
$ws.Range("J12").Value = 1.035407253819927
$ws.Range("K12").Value = 1.042469685822518
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.050368014429351
$ws.Range("N12").Value = 1.036877650883223

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028759239063825
$ws.Range("D13").Value = 1.03890912942164
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.046843725146626
$ws.Range("I13").Value = 1.037115465910658
$ws.Range("J13").Value = 1.035440097584511
$ws.Range("K13").Value = 1.042500702526131
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.050406004286951
$ws.Range("N13").Value = 1.03691054128972

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028928314201463
$ws.Range("D14").Value = 1.039042732190757
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.046999987370687
$ws.Range("I14").Value = 1.037154795843336
$ws.Range("J14").Value = 1.035547695206638
$ws.Range("K14").Value = 1.042602300685458
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.050530461552297
$ws.Range("N14").Value = 1.037018291712813

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029032487993983
$ws.Range("D15").Value = 1.039125042250723
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.047096266393975
$ws.Range("I15").Value = 1.037178985098143
$ws.Range("J15").Value = 1.035613976316187
$ws.Range("K15").Value = 1.042664875403979
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.050607128975844
$ws.Range("N15").Value = 1.037084666949139

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029638839785271
$ws.Range("D16").Value = 1.039604014483895
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.04765666200036
$ws.Range("I16").Value = 1.037319116801025
$ws.Range("J16").Value = 1.035999556937165
$ws.Range("K16").Value = 1.043028731577513
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.051053139810022
$ws.Range("N16").Value = 1.037470795138826

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030019201418757
$ws.Range("D17").Value = 1.039904365256335
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.048008191419935
$ws.Range("I17").Value = 1.037406435830224
$ws.Range("J17").Value = 1.036241241246196
$ws.Range("K17").Value = 1.043256655017258
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.051332710254662
$ws.Range("N17").Value = 1.037712822667297

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030241064553656
$ws.Range("D18").Value = 1.040079519906071
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.048213235164206
$ws.Range("I18").Value = 1.037457157681087
$ws.Range("J18").Value = 1.036382146472347
$ws.Range("K18").Value = 1.043389485526624
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.051495706479148
$ws.Range("N18").Value = 1.037853927995032

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030316714923165
$ws.Range("D19").Value = 1.040139237168693
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.048283150189304
$ws.Range("I19").Value = 1.037474416910293
$ws.Range("J19").Value = 1.036430180445774
$ws.Range("K19").Value = 1.04343475804395
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.051551271658442
$ws.Range("N19").Value = 1.03790203018221

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029978391733591
$ws.Range("D20").Value = 1.039872144060683
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.04797047536571
$ws.Range("I20").Value = 1.037397089028136
$ws.Range("J20").Value = 1.036215317540734
$ws.Range("K20").Value = 1.043232212710798
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.051302722483061
$ws.Range("N20").Value = 1.0376868621472

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028878525574359
$ws.Range("D21").Value = 1.039003390964885
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.046953971890779
$ws.Range("I21").Value = 1.037143223160212
$ws.Range("J21").Value = 1.035516013165993
$ws.Range("K21").Value = 1.042572387404189
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.050493815072416
$ws.Range("N21").Value = 1.036986564680036

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028187040304994
$ws.Range("D22").Value = 1.038456864267983
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.046314885984294
$ws.Range("I22").Value = 1.036981720354385
$ws.Range("J22").Value = 1.035075749961055
$ws.Range("K22").Value = 1.042156512820443
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.049984575750267
$ws.Range("N22").Value = 1.036545676250847

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028553605413886
$ws.Range("D23").Value = 1.038746617184703
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.046653674321
$ws.Range("I23").Value = 1.037067514774884
$ws.Range("J23").Value = 1.035309196860977
$ws.Range("K23").Value = 1.042377071882177
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.050254594205688
$ws.Range("N23").Value = 1.036779454672144

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029996831859305
$ws.Range("D24").Value = 1.039886703539855
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.047987517619203
$ws.Range("I24").Value = 1.037401313094115
$ws.Range("J24").Value = 1.036227031548631
$ws.Range("K24").Value = 1.043243257487308
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.051316272890227
$ws.Range("N24").Value = 1.037698592790332

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031671795967636
$ws.Range("D25").Value = 1.041208338423153
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.049535466015949
$ws.Range("I25").Value = 1.037780433156993
$ws.Range("J25").Value = 1.037289569968427
$ws.Range("K25").Value = 1.044243967279892
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.052545445210213
$ws.Range("N25").Value = 1.038762640136559
